$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L (2022 figures) mirroring the existing K column's
# formatting, then fill in the new values for each row.

# Row 2 (thin header separator row) - formatted like K2, no value.
$ws.Range("K2").Copy($ws.Range("L2"))

# Row 3 - year header, like K3 (2021) but 2022.
$ws.Range("K3").Copy($ws.Range("L3"))
$ws.Range("L3").Value = 2022

# Row 4
$ws.Range("K4").Copy($ws.Range("L4"))
$ws.Range("L4").Value = 370

# Row 5
$ws.Range("K5").Copy($ws.Range("L5"))
$ws.Range("L5").Value = 137

# Row 6
$ws.Range("K6").Copy($ws.Range("L6"))
$ws.Range("L6").Value = 314

# Row 7
$ws.Range("K7").Copy($ws.Range("L7"))
$ws.Range("L7").Value = 121

# Row 8 - this one picks up a thousands-separator number format (new style)
# rather than reusing K8's plain-General style.
$ws.Range("K8").Copy($ws.Range("L8"))
$ws.Range("L8").NumberFormat = "#,##0"
$ws.Range("L8").Value = 50

# Row 9
$ws.Range("K9").Copy($ws.Range("L9"))
$ws.Range("L9").Value = 16

# Move the active selection to L2 (was L5 before the edit).
$ws.Range("L2").Select()
